$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 61, pushing the existing rows 61-118 down to 63-120.
$ws.Rows.Item(61).Insert()
$ws.Rows.Item(61).Insert()

# New row 61 - "Primera" quality record for the latest week (2023-09-26)
$ws.Range("A61").Value = 7
$ws.Range("B61").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C61").Value = "Ñuble"
$ws.Range("D61").Value = 45195
$ws.Range("E61").Value = 16
$ws.Range("F61").Value = 100112044
$ws.Range("G61").Value = "Perejil"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 200
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = 1500
$ws.Range("N61").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O61").Value = "Región de Ñuble"
$ws.Range("P61").Value = 1500
$ws.Range("Q61").Value = 1
$ws.Range("R61").Value = "Hortaliza"

# New row 62 - "Segunda" quality record for the latest week (2023-09-26)
$ws.Range("A62").Value = 7
$ws.Range("B62").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C62").Value = "Ñuble"
$ws.Range("D62").Value = 45195
$ws.Range("E62").Value = 16
$ws.Range("F62").Value = 100112044
$ws.Range("G62").Value = "Perejil"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Segunda"
$ws.Range("J62").Value = 200
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = 1000
$ws.Range("N62").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O62").Value = "Región de Ñuble"
$ws.Range("P62").Value = 1000
$ws.Range("Q62").Value = 1
$ws.Range("R62").Value = "Hortaliza"

# Match the date number-format used by the rest of column D.
$ws.Range("D61").NumberFormat = $ws.Range("D63").NumberFormat
$ws.Range("D62").NumberFormat = $ws.Range("D63").NumberFormat
